$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metricas")

# Fill in the new "implementar Cola Dinamica" task data in row 8
$ws.Range("A8").Value = "implementar Cola Dinamica"
$ws.Range("B8").Value = 40
$ws.Range("C8").Value = 36
$ws.Range("D8").Value = 0.020833333333333332
$ws.Range("E8").Value = 0.11944444444444445
$ws.Range("F8").Value = 0.1451388888888889
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0

# Recalculate so that dependent formulas (G8, J8, B11, C11, D11, G11, H11,
# I11, J11, B12, B15, B16, B17, B18, B19, C19, B20, C20) pick up new values
$excel.CalculateFullRebuild()

# Update the active cell selection to match the edited workbook state
$ws.Range("C14").Select()

$wb.Save()
